$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(738, 'Coereba flaveola', 766, 2),
    @(245, 'Ramphocelus dimidiatus', 229, 174),
    @(513, 'Basileuterus sp', 500, 2),
    @(786, 'Coereba flaveola', 731, 3),
    @(737, 'Saltator maximus', 707, 14),
    @(761, 'Coereba flaveola', 723, 3),
    @(39, 'Hyloxalus', 140, 22),
    @(792, 'Saltator maximus', 682, 8),
    @(142, 'Leptopogon superciliaris', 114, 127),
    @(1049, 'Diasporus gularis', 795, 14),
    @(562, 'Call bird 2', 619, 81),
    @(526, 'Sing bird 1', 518, 52),
    @(114, 'Leptopogon superciliaris', 129, 4),
    @(129, 'Leptopogon superciliaris', 114, 4),
    @(385, 'Hyloxalus', 297, 6),
    @(389, 'Ramphocelus dimidiatus', 396, 1),
    @(518, 'Sing bird 1', 473, 4),
    @(473, 'Sing bird 1', 518, 4),
    @(113, 'Leptopogon superciliaris', 119, 7),
    @(464, 'Sing bird 1', 480, 3),
    @(119, 'Leptopogon superciliaris', 95, 3),
    @(381, 'Ramphocelus dimidiatus', 358, 9),
    @(358, 'Ramphocelus dimidiatus', 381, 9),
    @(601, 'Call bird 2', 771, 4),
    @(185, 'Ramphocelus dimidiatus', 358, 14),
    @(552, 'Call bird 2', 777, 1),
    @(574, 'Call bird 2', 588, 1),
    @(57, 'Leptopogon superciliaris', 18, 30),
    @(558, 'Hyloxalus', 171, 1),
    @(941, 'Diasporus gularis', 807, 69),
    @(566, 'Call bird 2', 758, 7),
    @(437, 'Sing bird 1', 469, 27),
    @(193, 'Ramphocelus dimidiatus', 325, 1),
    @(469, 'Sing bird 1', 459, 3),
    @(572, 'Hyloxalus', 568, 3),
    @(459, 'Sing bird 1', 469, 3),
    @(807, 'Diasporus gularis', 941, 69),
    @(8, 'Leptopogon superciliaris', 32, 46),
    @(318, 'Ramphocelus dimidiatus', 322, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
